$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New login-view translation rows, appended below the existing data in
# --- natural (unsorted) entry order. A subsequent Sort (by column A) will
# --- reorder rows 3:8 alphabetically by key, matching the final layout.

# Row 5: loginView.button.signIn
$ws.Cells.Item(5,1).Value = "loginView.button.signIn"
$ws.Cells.Item(5,2).Value = "Sign in"
$ws.Cells.Item(5,3).Value = "登入"
$ws.Cells.Item(5,4).Value = "登入"

# Row 6: loginView.desc.loginID
$ws.Cells.Item(6,1).Value = "loginView.desc.loginID"
$ws.Cells.Item(6,2).Value = "login ID"
$ws.Cells.Item(6,3).Value = "登入名"
$ws.Cells.Item(6,4).Value = "登入名"

# Row 7: loginView.desc.password
$ws.Cells.Item(7,1).Value = "loginView.desc.password"
$ws.Cells.Item(7,2).Value = "Password"
$ws.Cells.Item(7,3).Value = "密码"
$ws.Cells.Item(7,4).Value = "密碼"

# Row 8: loginView.desc.plzSignIn
$ws.Cells.Item(8,1).Value = "loginView.desc.plzSignIn"
$ws.Cells.Item(8,2).Value = "Please sign in"
$ws.Cells.Item(8,3).Value = "请登入帐户"
$ws.Cells.Item(8,4).Value = "請登入帳戶"

# Row 9: loginView.desc.rememberMe
$ws.Cells.Item(9,1).Value = "loginView.desc.rememberMe"
$ws.Cells.Item(9,2).Value = "Remember me"
$ws.Cells.Item(9,3).Value = "记住我"
$ws.Cells.Item(9,4).Value = "記住我"

# --- Widen column A to fit the longer i18next keys ---
$ws.Columns.Item(1).ColumnWidth = 31.28571428

# --- Sort the data rows (A3:D8) ascending by column A (the "key" column) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3"))
$ws.Sort.SetRange($ws.Range("A3:D8"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Move the active selection to C9, matching the post-edit cursor position ---
$ws.Range("C9").Select() | Out-Null
